# Replace the field `m:''.emptyText()` (a begin/instrText/end field-code
# triple, with an empty result) with a single literal run of text
# "{m:''.emptyText()}", wherever it appears among the document's
# headers/footers.
#
# NOTE: Field.Code.Start/.End/.Result use a coordinate space that does not
# line up with the owning header/footer Range in this runtime, so doing
# position arithmetic against a header/footer Range using numbers read
# from Field.Code/Field.Result lands edits in the wrong place. Field.Code
# .Text itself (the instruction text) reads back fine though, so that is
# all we rely on the Field object for.
#
# Approach, all done against the (correctly behaving) Range.Text /
# Range.Find / Range.InsertAfter API of the owning header/footer Range:
#   1. Find the field whose code contains "emptyText".
#   2. Split the Range's visible text on CR (paragraph marks) to figure
#      out, in plain text terms, which paragraph is the field's own
#      (empty, because its result is empty) paragraph, and capture the
#      text of the paragraph that follows it as a unique "anchor".
#   3. Delete the field (Field.Delete() correctly removes the
#      begin/instrText/end runs and leaves an empty paragraph behind).
#   4. Re-find that now-empty paragraph with Range.Find.Execute(anchor),
#      then collapse+step back one character to land inside it.
#   5. Build a *fresh* Range (header/footer .Range + SetRange) at that
#      position and InsertAfter() the literal text - reusing the Range
#      object that Find/Collapse/Move mutated can insert at a stale
#      location in this runtime, but a freshly obtained Range with the
#      same Start/End inserts correctly.

$d = $word.ActiveDocument
$CR = [char]13

function Fix-EmptyTextFields($rng) {
    $didChange = $false
    $guard = 0
    while ($rng.Fields.Count -gt 0 -and $guard -lt 50) {
        $guard = $guard + 1

        $fld = $null
        $fieldCount = $rng.Fields.Count
        $k = 1
        while ($k -le $fieldCount -and $fld -eq $null) {
            $candidate = $rng.Fields.Item($k)
            if ($candidate.Code.Text -like "*emptyText*") {
                $fld = $candidate
            }
            $k = $k + 1
        }
        if ($fld -eq $null) {
            break
        }

        # Work out, from the plain visible text, which paragraph is the
        # field's own (empty-result) paragraph, and what text follows it.
        $fullText = $rng.Text
        $parts = $fullText.Split($CR)
        $n = $parts.Length
        $lastIdx = $n - 1   # trailing split artifact when text ends on a CR

        $fieldParaIdx = -1
        $p = 0
        while ($p -lt $lastIdx -and $fieldParaIdx -lt 0) {
            if ($parts[$p].Length -eq 0) {
                $fieldParaIdx = $p
            }
            $p = $p + 1
        }

        $anchor = ""
        $haveAnchor = $false
        if ($fieldParaIdx -ge 0 -and ($fieldParaIdx + 1) -le $lastIdx) {
            $anchor = $parts[$fieldParaIdx + 1]
            if ($anchor.Length -gt 0) {
                $haveAnchor = $true
            }
        }

        $fld.Delete()

        $newText = "{m:''.emptyText()}"

        if ($haveAnchor) {
            $findRange = $rng.Duplicate
            $found = $findRange.Find.Execute($anchor, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
            if ($found) {
                $findRange.Collapse(1)
                $findRange.Move(1, -1)
                $insertRange = $rng.Duplicate
                $insertRange.SetRange($findRange.Start, $findRange.End)
                $insertRange.InsertAfter($newText)
                $didChange = $true
                continue
            }
        }

        # Fallback: no usable anchor (field paragraph was the last
        # paragraph in the story) - insert right at the end of the range.
        $insertRange = $rng.Duplicate
        $insertRange.SetRange($rng.End, $rng.End)
        $insertRange.InsertAfter($newText)
        $didChange = $true
    }
    return $didChange
}

foreach ($sec in $d.Sections) {
    foreach ($hdr in $sec.Headers) {
        if ($hdr.Range.Fields.Count -gt 0) {
            Fix-EmptyTextFields($hdr.Range) | Out-Null
        }
    }
    foreach ($ftr in $sec.Footers) {
        if ($ftr.Range.Fields.Count -gt 0) {
            Fix-EmptyTextFields($ftr.Range) | Out-Null
        }
    }
}

Write-Output "Done"
